# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to reflect the latest scraped counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "展览") {
        $ws.Range("F3").Value = 1131
        $ws.Range("F5").Value = 96
        $ws.Range("F8").Value = 11520
        $ws.Range("F9").Value = 4348
        $ws.Range("F14").Value = 2534
        $ws.Range("F18").Value = 2221
        $ws.Range("F19").Value = 179
        $ws.Range("F20").Value = 511
        $ws.Range("F21").Value = 11299
        $ws.Range("F22").Value = 11211
        $ws.Range("F24").Value = 41
        $ws.Range("F27").Value = 42
    }
    elseif ($sheetName -eq "全部类型") {
        $ws.Range("F3").Value = 1131
        $ws.Range("F5").Value = 96
        $ws.Range("F8").Value = 11520
        $ws.Range("F9").Value = 4348
        $ws.Range("F14").Value = 2534
        $ws.Range("F19").Value = 2221
        $ws.Range("F20").Value = 179
        $ws.Range("F21").Value = 511
        $ws.Range("F22").Value = 11299
        $ws.Range("F23").Value = 11211
        $ws.Range("F25").Value = 41
        $ws.Range("F28").Value = 42
    }
}
